$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.804.44"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "3.417.97"
$ws.Range("E3").Value = "  +3.83%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "259.05"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").Value = "671.04"
$ws.Range("E6").Value = "  +7.99%  "
$ws.Range("E7").Value = "  +10.33%  "
$ws.Range("D8").Value = "0.472"
$ws.Range("E8").Value = "  +19.05%  "
$ws.Range("E9").Value = "  +21.85%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "3.413.87"
$ws.Range("E12").Value = "  +11.18%  "
$ws.Range("D13").Value = "43.23"
$ws.Range("E13").Value = "  +13.79%  "
$ws.Range("E14").Value = "  +12.54%  "
$ws.Range("D15").Value = "6.13"
$ws.Range("E15").Value = "  +11.90%  "
$ws.Range("D16").Value = "98.376.19"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "4.062.23"
$ws.Range("E17").Value = "  +4.02%  "
$ws.Range("E18").Value = "  +34.89%  "
$ws.Range("D19").Value = "3.421.06"
$ws.Range("E19").Value = "  +3.44%  "
$ws.Range("D20").Value = "17.39"
$ws.Range("E20").Value = "  +15.19%  "
$ws.Range("D21").Value = "537.14"
$ws.Range("E21").Value = "  +12.72%  "
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("D23").Value = "10.68"
$ws.Range("E23").Value = "  +15.13%  "
$ws.Range("D24").Value = "0.0000219"
$ws.Range("E24").Value = "  +8.00%  "
$ws.Range("D25").Value = "0.439"
$ws.Range("E25").Value = "  +51.57%  "
$ws.Range("E26").Value = "  +15.43%  "
$ws.Range("D27").Value = "103.30"
$ws.Range("E27").Value = "  +17.28%  "
$ws.Range("E28").Value = "  +8.53%  "
$ws.Range("D29").Value = "3.598.68"
$ws.Range("E29").Value = "  +3.85%  "
$ws.Range("D30").Value = "0.151"
$ws.Range("E30").Value = "  +15.83%  "
$ws.Range("D31").Value = "11.72"
$ws.Range("E31").Value = "  +20.48%  "
$ws.Range("D32").Value = "0.199"
$ws.Range("E32").Value = "  +7.20%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("E35").Value = "  +10.85%  "
$ws.Range("E36").Value = "  +24.27%  "
$ws.Range("E37").Value = "  +14.89%  "
$ws.Range("D38").Value = "7.94"
$ws.Range("E38").Value = "  +11.62%  "
$ws.Range("E39").Value = "  +10.45%  "
$ws.Range("D40").Value = "534.07"
$ws.Range("E40").Value = "  +8.44%  "
$ws.Range("D41").Value = "1.43"
$ws.Range("E41").Value = "  +15.96%  "
$ws.Range("D42").Value = "24.79"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "0.0439"
$ws.Range("E43").Value = "  +36.05%  "
$ws.Range("E44").Value = "  +3.90%  "
$ws.Range("E45").Value = "  +11.38%  "
$ws.Range("D46").Value = "0.858"
$ws.Range("E46").Value = "  +8.62%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "8.24"
$ws.Range("E48").Value = "  +19.04%  "
$ws.Range("D49").Value = "2.13"
$ws.Range("E49").Value = "  +13.15%  "
$ws.Range("D50").Value = "5.33"
$ws.Range("E50").Value = "  +16.54%  "
$ws.Range("D51").Value = "1.59"
$ws.Range("E51").Value = "  +17.84%  "
